$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.168.84"
$ws.Range("E2").Value = "  -2.35%  "

$ws.Range("D3").Value = "1.871.29"
$ws.Range("E3").Value = "  -1.69%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").Value = "'307.47"
$ws.Range("E5").Value = "  -1.75%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").Value = "'0.5161"
$ws.Range("E7").Value = "  +2.42%  "

$ws.Range("D8").Value = "'0.3756"
$ws.Range("E8").Value = "  -1.49%  "

$ws.Range("D9").Value = "'0.07164"
$ws.Range("E9").Value = "  -1.49%  "

$ws.Range("D10").Value = "'20.86"
$ws.Range("E10").Value = "  +0.27%  "

$ws.Range("D11").Value = "'0.8858"
$ws.Range("E11").Value = "  -2.47%  "

$ws.Range("D12").Value = "1.891.49"
$ws.Range("E12").Value = "  -2.09%  "

$ws.Range("D13").Value = "'0.07580"
$ws.Range("E13").Value = "  -0.91%  "

$ws.Range("D14").Value = "'5.336"
$ws.Range("E14").Value = "  -2.57%  "

$ws.Range("D15").Value = "'89.48"
$ws.Range("E15").Value = "  -1.97%  "

$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.33%  "

$ws.Range("D17").Value = "'0.000008552"
$ws.Range("E17").Value = "  -1.77%  "

$ws.Range("D18").Value = "'14.20"
$ws.Range("E18").Value = "  -2.18%  "

$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").Value = "27.208.77"
$ws.Range("E20").Value = "  -2.31%  "

$ws.Range("D21").Value = "'5.037"
$ws.Range("E21").Value = "  -2.36%  "

$ws.Range("D22").Value = "2.119.26"
$ws.Range("E22").Value = "  -1.86%  "

$ws.Range("E23").Value = "  -1.42%  "

$ws.Range("D24").Value = "'6.486"

$ws.Range("D25").Value = "'151.60"
$ws.Range("E25").Value = "  -1.67%  "

$ws.Range("D26").Value = "'1.850"
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.05"
$ws.Range("E27").Value = "  -1.69%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.176"
$ws.Range("E28").Value = "  -2.50%  "

$ws.Range("D29").Value = "'113.24"
$ws.Range("E29").Value = "  -1.67%  "

$ws.Range("D30").Value = "'4.761"
$ws.Range("E30").Value = "  -2.94%  "

$ws.Range("D31").Value = "'4.703"
$ws.Range("E31").Value = "  +1.42%  "

$ws.Range("D32").Value = "'0.09034"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").Value = "'0.05172"
$ws.Range("E33").Value = "  -1.47%  "

$ws.Range("D34").Value = "'3.094"
$ws.Range("E34").Value = "  -3.27%  "

$ws.Range("D35").Value = "'0.7582"
$ws.Range("E35").Value = "  -0.76%  "

$ws.Range("D36").Value = "'1.181"
$ws.Range("E36").Value = "  -4.00%  "

$ws.Range("E37").Value = "  -0.57%  "

$ws.Range("D38").Value = "'2.511"
$ws.Range("E38").Value = "  -0.92%  "

$ws.Range("D39").Value = "'3.036"
$ws.Range("E39").Value = "  +0.58%  "

$ws.Range("D40").Value = "'1.082"
$ws.Range("E40").Value = "  -1.21%  "

$ws.Range("D41").Value = "'0.5426"
$ws.Range("E41").Value = "  -2.37%  "

$ws.Range("D42").Value = "'6.684"
$ws.Range("E42").Value = "  -4.03%  "

$ws.Range("D43").Value = "'115.11"
$ws.Range("E43").Value = "  +3.44%  "

$ws.Range("D44").Value = "'8.566"
$ws.Range("E44").Value = "  +1.13%  "

$ws.Range("D45").Value = "'0.1492"
$ws.Range("E45").Value = "  -1.23%  "

$ws.Range("D46").Value = "'0.4695"
$ws.Range("E46").Value = "  -2.00%  "

$ws.Range("D47").Value = "'10.17"
$ws.Range("E47").Value = "  -3.95%  "

$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("D49").Value = "'1.579"
$ws.Range("E49").Value = "  -3.11%  "

$ws.Range("D50").Value = "'64.98"
$ws.Range("E50").Value = "  -3.44%  "

$ws.Range("D51").Value = "'36.53"
